$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.320.51'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '2.647.51'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.54%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("D9").Value = '2.646.34'
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("E10").Value = '  -2.26%  '
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("D14").Value = '3.135.92'
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000187'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.10%  '
$ws.Range("D16").Value = '72.260.71'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.01'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.05%  '
$ws.Range("D18").Value = '2.605.29'
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '370.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.23%  '
$ws.Range("D28").Value = '2.788.13'
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("D30").Value = '0.0₃0961'
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '501.59'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.48%  '
$ws.Range("E33").Value = '  -1.82%  '
$ws.Range("E34").Value = '  -0.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.68'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.117'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("E40").Value = '  -2.92%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.57'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.92'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.330'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '154.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.552'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0749'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.29%  '
